# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback datetime
# columns for the f1607909-...-md file row across the Overview, zh-cn and
# de-de sheets, reflecting a newer report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to f1607909-33e9-459d-84e0-46fd80efbca0.md
$wsOverview.Range("G3").Value = "2016-08-19 04:43:18"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to f1607909-33e9-459d-84e0-46fd80efbca0.md
$wsZhCn.Range("H3").Value = "2016-08-19 04:43:14"
$wsZhCn.Range("K3").Value = "2016-08-19 04:43:30"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to f1607909-33e9-459d-84e0-46fd80efbca0.md
$wsDeDe.Range("H3").Value = "2016-08-19 04:43:18"
$wsDeDe.Range("K3").Value = "2016-08-19 04:43:37"
